# Update "activity till excel form" - refresh stats (runs, balls, fours, sixes)
# for Ishan Kishan / Mumbai Indians sheet. The underlying data rows were
# reshuffled/updated, so capture the current values first, then write back
# the new values per row for columns C (runs), D (balls), E (fours), F (sixes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only rows whose runs/balls/fours/sixes actually changed are touched;
# rows 7, 10, 11 and 12 stay exactly as they were.
$newValues = @{
    2  = @("68", "37", "6", "5")
    3  = @("7",  "7",  "1", "0")
    4  = @("72", "47", "8", "3")
    5  = @("33", "30", "1", "2")
    6  = @("25", "19", "3", "1")
    8  = @("55", "30", "4", "3")
    9  = @("0",  "1",  "0", "0")
    13 = @("99", "58", "2", "9")
    14 = @("28", "32", "1", "1")
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
    $ws.Cells.Item($row, 6).Value = $vals[3]
}
